$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.748.84"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.678.57"
$ws.Range("E3").Value = "  -1.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'600.77"
$ws.Range("E5").Value = "  -1.36%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'156.77"
$ws.Range("E6").Value = "  -0.87%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +4.16%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -1.97%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.74%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.05%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'29.30"
$ws.Range("E13").Value = "  -3.86%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -3.04%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.159.97"
$ws.Range("E15").Value = "  -1.11%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "65.618.82"
$ws.Range("E16").Value = "  -0.27%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.681.00"
$ws.Range("E17").Value = "  -0.91%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'12.90"
$ws.Range("E18").Value = "  +1.36%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -2.06%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.88%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'352.99"
$ws.Range("E21").Value = "  -1.89%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.00%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'69.88"
$ws.Range("E23").Value = "  -1.47%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +5.56%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "  -2.33%  "

# Row 26 - SuiNetwork
$ws.Range("E26").Value = "  -0.44%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -2.00%  "

# Row 28 - Fetch.AI
$ws.Range("D28").Value = "'1.60"
$ws.Range("E28").Value = "  -6.13%  "

# Row 29 - Aptos
$ws.Range("D29").Value = "'8.06"
$ws.Range("E29").Value = "  -4.44%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.02%  "

# Row 31 - Bittensor
$ws.Range("D31").Value = "'533.20"
$ws.Range("E31").Value = "  -1.89%  "

# Row 32 - PancakeSwap
$ws.Range("D32").Value = "'2.12"
$ws.Range("E32").Value = "  -3.97%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  +1.97%  "

# Row 35 - RenderToken
$ws.Range("E35").Value = "  -3.30%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("D36").Value = "'0.423"
$ws.Range("E36").Value = "  -2.66%  "

# Row 37 - EthereumClassic
$ws.Range("E37").Value = "  -1.54%  "

# Row 38 - was Monero, now FirstDigitalUSD
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.01%  "

# Row 39 - was FirstDigitalUSD, now Monero
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'159.30"
$ws.Range("E39").Value = "  -2.56%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -2.64%  "

# Row 42 - Aave
$ws.Range("D42").Value = "'163.92"
$ws.Range("E42").Value = "  -5.41%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  -1.62%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  +2.26%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "'0.0610"
$ws.Range("E45").Value = "  -1.23%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "'22.83"
$ws.Range("E46").Value = "  -3.03%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  -3.19%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  -2.44%  "

# Row 49 - BabyDogeCoin (contains subscript six, U+2086)
$ws.Range("D49").Value = "0.0{0}0264" -f [char]0x2086
$ws.Range("E49").Value = "  +15.10%  "

# Row 50 - EnergySwap
$ws.Range("E50").Value = "  -3.82%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +0.50%  "
